$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Actualización desde MV -datos-": append the new daily rows that follow
# the existing data (rows 2..251, ending at 07-09-2021) with dates
# 08-09-2021 .. 15-09-2021, starting at row 252.
$dates = @("08-09-2021", "09-09-2021", "10-09-2021", "11-09-2021", "12-09-2021", "13-09-2021", "14-09-2021", "15-09-2021")
$montos = @(6460, 6381, 6351, 6351, 6351, 6362, 6408, 6416)
$stocks = @(13, 13, 13, 13, 13, 13, 13, 13)

$startRow = 252
$endRow = $startRow + $dates.Count - 1

# Column A holds dates written as plain text (e.g. "08-09-2021"), matching
# the existing rows. Force text formatting first so Excel's automatic date
# recognition doesn't turn these into date serial numbers, then restore the
# default "Normal" style so the new cells look just like the rest of the
# sheet (no explicit style index).
$rngA = $ws.Range("A$startRow`:A$endRow")
$rngA.NumberFormat = "@"

for ($i = 0; $i -lt $dates.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $montos[$i]
    $ws.Cells.Item($r, 3).Value = $stocks[$i]
}

$rngA.Style = "Normal"
